# Commit: "updated end of day for all 4ml horz runs of beads"
# Appends five new rows (19-23) of 4ml horizontal bead-run data to Sheet1,
# logging the end-of-day timestamps/readings for those runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 19: D20151103T202347
$ws.Cells.Item(19, 1).Value = "D20151103T202347"
$ws.Cells.Item(19, 2).Value = 4
$ws.Cells.Item(19, 3).Value = "H"
$ws.Cells.Item(19, 4).Value = 560.4
$ws.Cells.Item(19, 5).Value = 3.3262999999999998
$ws.Cells.Item(19, 6).Value = 958
$ws.Cells.Item(19, 7).Value = 144.69999999999999
$ws.Cells.Item(19, 8).Value = 1864
$ws.Cells.Item(19, 9).Value = "9um beads, use all signals"

# Row 20: D20151103T204202
$ws.Cells.Item(20, 1).Value = "D20151103T204202"
$ws.Cells.Item(20, 2).Value = 4
$ws.Cells.Item(20, 3).Value = "H"
$ws.Cells.Item(20, 4).Value = 575.70000000000005
$ws.Cells.Item(20, 5).Value = 3.3142999999999998
$ws.Cells.Item(20, 6).Value = 958
$ws.Cells.Item(20, 7).Value = 148.07
$ws.Cells.Item(20, 8).Value = 1908
$ws.Cells.Item(20, 9).Value = "9um beads, use all signals"

# Row 21: D20151103T210016
$ws.Cells.Item(21, 1).Value = "D20151103T210016"
$ws.Cells.Item(21, 2).Value = 4
$ws.Cells.Item(21, 3).Value = "H"
$ws.Cells.Item(21, 4).Value = 599.29999999999995
$ws.Cells.Item(21, 5).Value = 3.2921
$ws.Cells.Item(21, 6).Value = 958
$ws.Cells.Item(21, 7).Value = 153.1
$ws.Cells.Item(21, 8).Value = 1973
$ws.Cells.Item(21, 9).Value = "9um beads, use all signals"

# Row 22: D20151103T211831
$ws.Cells.Item(22, 1).Value = "D20151103T211831"
$ws.Cells.Item(22, 2).Value = 4
$ws.Cells.Item(22, 3).Value = "H"
$ws.Cells.Item(22, 4).Value = 557.6
$ws.Cells.Item(22, 5).Value = 3.3338000000000001
$ws.Cells.Item(22, 6).Value = 958
$ws.Cells.Item(22, 7).Value = 143.5
$ws.Cells.Item(22, 8).Value = 1859
$ws.Cells.Item(22, 9).Value = "9um beads, use all signals"

# Row 23: D20151103T214021 - only Filename/Volume/HorzOrVert recorded so far
$ws.Cells.Item(23, 1).Value = "D20151103T214021"
$ws.Cells.Item(23, 2).Value = 3
$ws.Cells.Item(23, 3).Value = "H"

# Move the selection to reflect where the author left off editing
$ws.Range("H23").Select()
